# Update the "datetimeFigureOut" date field shown on the slide master and
# every slide layout from 22/11/2021 to 25/11/2021.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "22/11/2021") {
                $tr.Text = "25/11/2021"
            }
        }
    }
}

# Slide master's own Date Placeholder shape.
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout has its own Date Placeholder shape.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
